$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$letters = @("A", "B", "C")
$row = 1
foreach ($letter in $letters) {
    for ($n = 1; $n -le 11; $n++) {
        $ws.Cells.Item($row, 1).Value = $letter
        $ws.Cells.Item($row, 2).Value = $n
        $row++
    }
}

$ws.Range("A23:A33").Select()
